# Apply weekly "Fruta / hortaliza" update: insert two new data rows (new
# rows 415 and 416) into the "Vega Monumental Concepción - Naranja" sheet,
# pushing the previously-existing rows 415-449 down to 417-451.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows right before the current row 415, shifting
# everything from row 415 down (old 415 -> new 417, ... old 449 -> new 451).
$ws.Rows("415:416").Insert()

# --- New row 415 -----------------------------------------------------
$ws.Range("A415").Value = 11
$ws.Range("B415").Value = "Vega Monumental Concepción"
$ws.Range("C415").Value = "Bíobío"
$ws.Range("D415").Value = "2023-03-28"
$ws.Range("E415").Value = 8
$ws.Range("F415").Value = "Fruta"
$ws.Range("G415").Value = 100102
$ws.Range("H415").Value = "Cítricos"
$ws.Range("I415").Value = 100102005
$ws.Range("J415").Value = "Naranja"
$ws.Range("K415").Value = "Valencia"
$ws.Range("L415").Value = "Primera"
$ws.Range("M415").Value = 200
$ws.Range("N415").Value = 14000
$ws.Range("O415").Value = 15000
$ws.Range("P415").Value = 14500
$ws.Range("Q415").Value = "$/caja 15 kilos empedrada"
$ws.Range("R415").Value = "Región de O'Higgins"
$ws.Range("S415").Value = 967
$ws.Range("T415").Value = 15

# --- New row 416 -----------------------------------------------------
$ws.Range("A416").Value = 11
$ws.Range("B416").Value = "Vega Monumental Concepción"
$ws.Range("C416").Value = "Bíobío"
$ws.Range("D416").Value = "2023-03-28"
$ws.Range("E416").Value = 8
$ws.Range("F416").Value = "Fruta"
$ws.Range("G416").Value = 100102
$ws.Range("H416").Value = "Cítricos"
$ws.Range("I416").Value = 100102005
$ws.Range("J416").Value = "Naranja"
$ws.Range("K416").Value = "Valencia"
$ws.Range("L416").Value = "Segunda"
$ws.Range("M416").Value = 100
$ws.Range("N416").Value = 12000
$ws.Range("O416").Value = 12000
$ws.Range("P416").Value = 12000
$ws.Range("Q416").Value = "$/caja 15 kilos empedrada"
$ws.Range("R416").Value = "Región de O'Higgins"
$ws.Range("S416").Value = 800
$ws.Range("T416").Value = 15
